$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column before column A, shifting A:E to B:F.
$ws.Columns("A").Insert()

# New column A: header left blank, data rows get a 0-based index
# (used for the new "total heatmap" coloring), styled like the other
# header cells (bold, centered, bordered, same as row 1).
for ($i = 0; $i -le 15; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $i
}

$ws.Range("B1").Copy()
$ws.Range("A2:A17").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A1:F17").Select()
